$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''305.30'
$ws.Cells.Item(2, 5).Value = '''5.89%'
$ws.Cells.Item(3, 4).Value = '''35.37'
$ws.Cells.Item(3, 5).Value = '''13.82%'
$ws.Cells.Item(4, 4).Value = '''5.142'
$ws.Cells.Item(4, 5).Value = '''4.29%'
$ws.Cells.Item(5, 4).Value = '''0.07780'
$ws.Cells.Item(5, 5).Value = '''5.88%'
$ws.Cells.Item(6, 4).Value = '''2.411'
$ws.Cells.Item(6, 5).Value = '''6.99%'
$ws.Cells.Item(7, 4).Value = '''8.021'
$ws.Cells.Item(7, 5).Value = '''3.79%'
$ws.Cells.Item(8, 2).Value = 'GateToken'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(8, 4).Value = '''3.942'
$ws.Cells.Item(8, 5).Value = '''5.58%'
$ws.Cells.Item(9, 2).Value = 'MXToken'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(9, 4).Value = '''0.9224'
$ws.Cells.Item(9, 5).Value = '''1.87%'
$ws.Cells.Item(10, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(10, 4).Value = '''0.09983'
$ws.Cells.Item(10, 5).Value = '''13.92%'
$ws.Cells.Item(11, 2).Value = 'WazirX'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(11, 4).Value = '''0.1803'
$ws.Cells.Item(11, 5).Value = '''7.18%'
$ws.Cells.Item(12, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(12, 4).Value = '''0.08648'
$ws.Cells.Item(12, 5).Value = '''5.09%'
$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(13, 4).Value = '''0.03320'
$ws.Cells.Item(13, 5).Value = '''6.79%'
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14, 4).Value = '''0.09895'
$ws.Cells.Item(14, 5).Value = '''-0.55%'
$ws.Cells.Item(15, 2).Value = 'BitForexToken'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(15, 4).Value = '''0.001497'
$ws.Cells.Item(15, 5).Value = '''-0.07%'
$ws.Cells.Item(16, 2).Value = 'TigerCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(16, 4).Value = '''0.005665'
$ws.Cells.Item(16, 5).Value = '''-2.34%'
$ws.Cells.Item(17, 2).Value = 'LEO'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(17, 4).Value = '''3.470'
$ws.Cells.Item(17, 5).Value = '''-0.49%'
$ws.Cells.Item(18, 4).Value = '''2.141'
$ws.Cells.Item(18, 5).Value = '''3.80%'
$ws.Cells.Item(19, 4).Value = '''0.3364'
$ws.Cells.Item(19, 5).Value = '''1.04%'
$ws.Cells.Item(20, 5).Value = '''0.23%'
$ws.Cells.Item(21, 4).Value = '''4.308'
$ws.Cells.Item(21, 5).Value = '''3.44%'
$ws.Cells.Item(22, 4).Value = '''0.2385'
$ws.Cells.Item(22, 5).Value = '''12.30%'
$ws.Cells.Item(23, 4).Value = '''0.04570'
$ws.Cells.Item(23, 5).Value = '''0.78%'
$ws.Cells.Item(24, 5).Value = '''0.54%'
$ws.Cells.Item(25, 4).Value = '''0.004457'
$ws.Cells.Item(25, 5).Value = '''7.48%'
$ws.Cells.Item(26, 4).Value = '''0.0001299'
$ws.Cells.Item(26, 5).Value = '''-0.12%'
$ws.Cells.Item(27, 5).Value = '''8.74%'
$ws.Cells.Item(39, 4).Value = '''0.01791'
$ws.Cells.Item(39, 5).Value = '''13.68%'
$ws.Cells.Item(40, 4).Value = '''0.04754'
$ws.Cells.Item(40, 5).Value = '''6.92%'
$ws.Cells.Item(41, 4).Value = '''0.007740'
$ws.Cells.Item(41, 5).Value = '''5.81%'
$ws.Cells.Item(42, 5).Value = '''6.75%'
$ws.Cells.Item(43, 4).Value = '''0.007093'
$ws.Cells.Item(43, 5).Value = '''-25.93%'
$ws.Cells.Item(44, 4).Value = '''0.002104'
$ws.Cells.Item(44, 5).Value = '''-5.67%'
$ws.Cells.Item(45, 4).Value = '''0.009534'
$ws.Cells.Item(45, 5).Value = '''13.09%'
$ws.Cells.Item(46, 4).Value = '''0.00006114'
$ws.Cells.Item(46, 5).Value = '''-0.21%'
$ws.Cells.Item(47, 5).Value = '''-0.25%'
$ws.Cells.Item(48, 5).Value = '''25.04%'
$ws.Cells.Item(49, 5).Value = '''-0.25%'
$ws.Cells.Item(50, 4).Value = '''0.00002098'
$ws.Cells.Item(50, 5).Value = '''-0.25%'
$ws.Cells.Item(51, 4).Value = '''0.0001998'
$ws.Cells.Item(51, 5).Value = '''-0.25%'
